# R22 UAT2 - Regression: duplicate the amendmentLockerJoint sheet, place the
# duplicate ("UAT2 Regression") before it, and update the original sheet's
# KEY value in A2.

$wb = $excel.ActiveWorkbook

# Locate the original worksheet and copy it in place before itself -- this
# creates a second worksheet ("amendmentLockerJoint (2)") positioned first.
$orig = $wb.Worksheets.Item("amendmentLockerJoint")
$orig.Copy($orig)

# NOTE: after Copy() the $orig handle tracks the newly-created copy (it
# becomes the active sheet), so re-resolve both sheets by position/name
# rather than reusing $orig.
$newSheet = $wb.Worksheets.Item(1)
$newSheet.Name = "UAT2 Regression"

$amendment = $wb.Worksheets.Item("amendmentLockerJoint")

# Update the KEY cell on the original sheet only; the duplicate keeps the
# old value.
$amendment.Range("A2").Value = "OR.0003.0031"

# Keep the original sheet as the active/selected tab, matching the source
# workbook's prior state.
$amendment.Activate()
